$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating "2022-Q3" (this gives us
#    an exact copy of its formatting / styles) and inserting it right before
#    the existing "2022-Q3" tab, then rename it.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3) | Out-Null

$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# The copied sheet has 4 rows of data (3 funds); the 2022-Q4 snapshot only
# has 2 funds, so drop the extra row.
$newSheet.Rows.Item(4).Delete()

# Force the text-like columns (fund code / name / scale / position / ratio /
# market value) to stay text so leading zeros (e.g. "006348") and the exact
# original formatting ("11.81") are preserved instead of being turned into
# numbers.
$newSheet.Range("B2:G3").NumberFormat = "@"

# Row 2: 006348 / 银华盛利混合A
$newSheet.Range("B2").Value() = "006348"
$newSheet.Range("C2").Value() = "银华盛利混合A"
$newSheet.Range("D2").Value() = "11.81"
$newSheet.Range("E2").Value() = "86.13"
$newSheet.Range("F2").Value() = "2.78"
$newSheet.Range("G2").Value() = "0.3283"
$newSheet.Range("H2").Value() = 8

# Row 3: 015684 / 银华盛利混合C
$newSheet.Range("B3").Value() = "015684"
$newSheet.Range("C3").Value() = "银华盛利混合C"
$newSheet.Range("D3").Value() = "1.68"
$newSheet.Range("E3").Value() = "86.13"
$newSheet.Range("F3").Value() = "2.78"
$newSheet.Range("G3").Value() = "0.0467"
$newSheet.Range("H3").Value() = 8

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new row for 2022-Q4 above the
#    existing 2022-Q3 row, shifting the rest down, and renumber the index
#    column.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# The index column (A) uses a bold / centered / bordered style; copy it from
# the row below (which still has the original style) onto the new row.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$total.Range("A2").Value() = 0
$total.Range("B2").Value() = "2022-Q4"
$total.Range("C2").Value() = 2
$total.Range("D2").Value() = 0.38

$total.Range("A3").Value() = 1
$total.Range("A4").Value() = 2

# Restore the tab that was originally selected in the source workbook.
$wb.Worksheets.Item("2021-Q3").Activate()
